$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the now-empty C4 cell (it referenced an empty shared string, removed by this edit)
$ws.Range("C4").Value = $null

# New row 5 - "Euro kuru" form submission
$ws.Range("A5").Value = 43158.08285303263
$ws.Range("A5").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("C5").Value = "Euro kuru"
$ws.Range("I5").Value = "http://euro.tlkur.com"

# New row 6 - "Euro kuru" form submission
$ws.Range("A6").Value = 43158.09071728615
$ws.Range("A6").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("C6").Value = "Euro kuru"
$ws.Range("I6").Value = "http://euro.tlkur.com/"
